$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 0.1876208897485493
    "C2" = 0.5570599613152805
    "J2" = 0.007736943907156673
    "P2" = 0.1411992263056093
    "S2" = 0.1063829787234043
    "B3" = 0.01628664495114007
    "C3" = 0.03908794788273615
    "J3" = 0.03583061889250815
    "P3" = 0.6677524429967426
    "S3" = 0.241042345276873
    "J4" = 0.0759493670886076
    "O4" = 0.02531645569620253
    "P4" = 0.6075949367088608
    "S4" = 0.2911392405063291
    "B6" = 0.04680851063829787
    "D6" = 0.00425531914893617
    "F6" = 0.09148936170212765
    "J6" = 0.2319148936170213
    "O6" = 0.0148936170212766
    "Q6" = 0.1574468085106383
    "R6" = 0.06808510638297872
    "S6" = 0.3851063829787234
    "B7" = 0.1093023255813954
    "D7" = 0.02093023255813953
    "F7" = 0.06046511627906977
    "J7" = 0.1046511627906977
    "O7" = 0.02558139534883721
    "Q7" = 0.1581395348837209
    "R7" = 0.08372093023255814
    "S7" = 0.4372093023255814
    "B8" = 0.07861936720997123
    "D8" = 0.01629913710450623
    "F8" = 0.07094918504314478
    "J8" = 0.1112176414189837
    "O8" = 0.01821668264621285
    "Q8" = 0.1658676893576222
    "R8" = 0.1045062320230105
    "S8" = 0.4343240651965484
    "B9" = 0.07731958762886598
    "D9" = 0.01288659793814433
    "F9" = 0.06958762886597938
    "J9" = 0.09793814432989691
    "O9" = 0.01288659793814433
    "Q9" = 0.2061855670103093
    "R9" = 0.08247422680412371
    "S9" = 0.4407216494845361
    "B10" = 0.08941176470588236
    "D10" = 0.01843137254901961
    "E10" = 0.001568627450980392
    "F10" = 0.06784313725490196
    "J10" = 0.1113725490196078
    "O10" = 0.01215686274509804
    "Q10" = 0.2109803921568627
    "R10" = 0.08313725490196078
    "S10" = 0.4050980392156863
    "G11" = 0.1341991341991342
    "J11" = 0.1024531024531024
    "K11" = 0.2049062049062049
    "L11" = 0.5425685425685426
    "S11" = 0.01587301587301587
    "G12" = 0.7020725388601037
    "J12" = 0.227979274611399
    "K12" = 0.01036269430051814
    "L12" = 0.0155440414507772
    "S12" = 0.04404145077720207
    "G13" = 0.6893203883495146
    "J13" = 0.2815533980582524
    "S13" = 0.02912621359223301
    "F15" = 0.02669902912621359
    "H15" = 0.1820388349514563
    "I15" = 0.07281553398058252
    "J15" = 0.3106796116504854
    "K15" = 0.0558252427184466
    "M15" = 0.01213592233009709
    "N15" = 0.002427184466019417
    "O15" = 0.04854368932038835
    "S15" = 0.2888349514563107
    "F16" = 0.009375
    "H16" = 0.19375
    "I16" = 0.053125
    "J16" = 0.390625
    "K16" = 0.115625
    "M16" = 0.028125
    "O16" = 0.053125
    "S16" = 0.15625
    "F17" = 0.00968783638320775
    "H17" = 0.2055974165769645
    "I17" = 0.09149623250807319
    "J17" = 0.4348762109795479
    "K17" = 0.08073196986006459
    "M17" = 0.01506996770721206
    "O17" = 0.05920344456404737
    "S17" = 0.1033369214208827
    "F18" = 0.02369668246445497
    "H18" = 0.2061611374407583
    "I18" = 0.07582938388625593
    "J18" = 0.3981042654028436
    "K18" = 0.1113744075829384
    "M18" = 0.004739336492890996
    "N18" = 0.002369668246445498
    "O18" = 0.07109004739336493
    "S18" = 0.1066350710900474
    "F19" = 0.0161059413027917
    "H19" = 0.2269148174659986
    "I19" = 0.08052970651395848
    "J19" = 0.3443092340730136
    "K19" = 0.1309949892627058
    "M19" = 0.02755905511811024
    "O19" = 0.05619183965640658
    "S19" = 0.117394416607015
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
